$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text columns (Coin name, Link) - no special formatting needed
$textValues = @{
    "B9" = "WazirX"
    "C9" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "B10" = "MandalaExchangeToken"
    "C10" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "B11" = "BitrueCoin"
    "C11" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "B12" = "BitMartToken"
    "C12" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "B13" = "BitForexToken"
    "C13" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "B14" = "CoinExToken"
    "C14" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "B15" = "One"
    "C15" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "B41" = "BKEXToken"
    "C41" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "B42" = "KickToken"
    "C42" = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
}

foreach ($key in $textValues.Keys) {
    $ws.Range($key).Value = $textValues[$key]
}

# Numeric-looking text columns (Price, Volume) - force text format so the
# stored value keeps its exact string representation (e.g. trailing zeros, %)
$numericTextValues = @{
    "D2" = "244.36"
    "E2" = "-1.11%"
    "E3" = "3.04%"
    "D4" = "5.117"
    "E4" = "0.20%"
    "D5" = "0.05666"
    "E5" = "0.79%"
    "D6" = "6.468"
    "E6" = "-0.61%"
    "D7" = "0.8229"
    "E7" = "1.30%"
    "D8" = "0.8448"
    "E8" = "-0.21%"
    "D9" = "0.1325"
    "E9" = "-1.48%"
    "D10" = "0.06935"
    "E10" = "-0.31%"
    "D11" = "0.02885"
    "E11" = "2.85%"
    "D12" = "0.09387"
    "E12" = "-0.12%"
    "D13" = "0.001520"
    "E13" = "0.51%"
    "D14" = "0.04113"
    "E14" = "-12.18%"
    "D15" = "0.0006019"
    "E15" = "0.73%"
    "D16" = "0.006219"
    "E16" = "0.73%"
    "D17" = "3.512"
    "E17" = "-1.76%"
    "E18" = "-1.80%"
    "D19" = "2.227"
    "E19" = "5.13%"
    "D21" = "0.03167"
    "E21" = "0.77%"
    "D23" = "3.559"
    "E23" = "-5.54%"
    "E24" = "-0.01%"
    "D25" = "0.001219"
    "E26" = "-3.86%"
    "E27" = "2.12%"
    "E28" = "3.53%"
    "E40" = "0.08%"
    "D41" = "0.1053"
    "E41" = "-22.34%"
    "D42" = "0.003440"
    "E42" = "-43.57%"
    "D43" = "0.002290"
    "E43" = "-13.88%"
    "D44" = "0.009684"
    "E44" = "11.69%"
    "D45" = "0.00005317"
    "E45" = "0.47%"
    "E46" = "0.07%"
    "E47" = "-15.78%"
    "D48" = "0.002590"
    "E48" = "25.47%"
    "E49" = "0.07%"
    "E50" = "0.07%"
}

foreach ($key in $numericTextValues.Keys) {
    $cell = $ws.Range($key)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextValues[$key]
}
